# Generate Report for Handoff
# Adds two newly-handed-off files (5ed8429d-... and 7e509f3e-...) to the
# localization status report: one new row each on the "Overview" sheet and
# on each language sheet ("zh-cn", "de-de"), pushing the existing
# ".localization-config" summary row down to make room.

$wb = $excel.ActiveWorkbook

$file1 = "5ed8429d-15b8-40b0-900e-f8d2765d3450"
$file2 = "7e509f3e-3515-40fc-84df-8905deb8e596"
$hash1 = "0e058f4ddf1562cea42d4d85233521413ff0b478"
$hash2 = "aff9d82b38fe8eed489c8020d5203f6b4bab8ca8"

$statusReady = "Ready for handoff"
$statusNotLoc = "Not to be localized"
$cfgName = ".localization-config"

# ---------------------------------------------------------------------
# Sheet "Overview": columns A=File Name, B=zh-cn, C=de-de
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# push the ".localization-config" row from row 4 down to row 6
$ws.Cells.Item(6,1).Value = $ws.Cells.Item(4,1).Value()
$ws.Cells.Item(6,2).Value = $ws.Cells.Item(4,2).Value()
$ws.Cells.Item(6,3).Value = $ws.Cells.Item(4,3).Value()

# new row for file1
$ws.Cells.Item(4,1).Value = $file1 + ".md"
$ws.Cells.Item(4,2).Value = $statusReady
$ws.Cells.Item(4,3).Value = $statusReady

# new row for file2
$ws.Cells.Item(5,1).Value = $file2 + ".md"
$ws.Cells.Item(5,2).Value = $statusReady
$ws.Cells.Item(5,3).Value = $statusReady

# rebuild hyperlinks (this engine can only append hyperlinks reliably, so
# clear the sheet's hyperlinks and recreate all of them in final order)
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/813411781d43403d6b153b57bceef52504c58b28/e2e/b78f0758-4fa2-4641-b5e7-7bb1f6d160c8.md", "", "", "b78f0758-4fa2-4641-b5e7-7bb1f6d160c8.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/efb688dfdd2528f758c385524c64ff5590bae1c0/e2e/56c59f66-ca12-44db-805b-c48245e4767b.md", "", "", "56c59f66-ca12-44db-805b-c48245e4767b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b/e2e/" + $file1 + ".md", "", "", $file1 + ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c/e2e/" + $file2 + ".md", "", "", $file2 + ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/813411781d43403d6b153b57bceef52504c58b28/.localization-config", "", "", $cfgName) | Out-Null

# ---------------------------------------------------------------------
# Language sheets "zh-cn" / "de-de"
# columns: A Source File Name, B Status, C Latest Handoff File,
#          D Latest Handoff Datetime, E Latest Target File,
#          F Latest Handback File, G Latest Handback DateTime,
#          H Handoff Reason, I Dependency From
# ---------------------------------------------------------------------
$langs = @(
    @{ Name = "zh-cn"; HandoffDT = "2016-03-10 04:44:04"; AUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b/e2e/" + $file1 + ".md"; AUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c/e2e/" + $file2 + ".md"; CUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c2d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $file1 + "." + $hash1 + ".zh-cn.xlf"; CUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c2d3e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $file2 + "." + $hash2 + ".zh-cn.xlf"; A2Url = "https://github.com/OpenLocalizationTest/oltest/blob/813411781d43403d6b153b57bceef52504c58b28/e2e/b78f0758-4fa2-4641-b5e7-7bb1f6d160c8.md"; C2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dddca6e62abe75a931346cdab094f0c7ee5b882c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b78f0758-4fa2-4641-b5e7-7bb1f6d160c8.cc3665e5c12c629e1fd08d7de85a74ca237eb321.zh-cn.xlf"; E2Url = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/d4e67e2fac3280292cffb7a3f632209f3e0719d9/e2e/b78f0758-4fa2-4641-b5e7-7bb1f6d160c8.md"; F2Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3f3b28184658943ccdc7213f1d9cd796360101a4/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b78f0758-4fa2-4641-b5e7-7bb1f6d160c8.cc3665e5c12c629e1fd08d7de85a74ca237eb321.zh-cn.xlf"; A3Url = "https://github.com/OpenLocalizationTest/oltest/blob/efb688dfdd2528f758c385524c64ff5590bae1c0/e2e/56c59f66-ca12-44db-805b-c48245e4767b.md"; C3Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71236da602bc62ed92a086656fdb7a867fb6ef60/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/56c59f66-ca12-44db-805b-c48245e4767b.fbc4bafe5dd05374473974c13a876717daac1cad.zh-cn.xlf" },
    @{ Name = "de-de"; HandoffDT = "2016-03-10 04:44:08"; AUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b/e2e/" + $file1 + ".md"; AUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c/e2e/" + $file2 + ".md"; CUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c2d3e4f/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $file1 + "." + $hash1 + ".de-de.xlf"; CUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c2d3e4f5a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $file2 + "." + $hash2 + ".de-de.xlf"; A2Url = "https://github.com/OpenLocalizationTest/oltest/blob/813411781d43403d6b153b57bceef52504c58b28/e2e/b78f0758-4fa2-4641-b5e7-7bb1f6d160c8.md"; C2Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/60b81f799f0dfeb11bbd6910b5f343ab9695a2ec/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b78f0758-4fa2-4641-b5e7-7bb1f6d160c8.cc3665e5c12c629e1fd08d7de85a74ca237eb321.de-de.xlf"; E2Url = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8e717367b7ac92b0e7f2f30a5ca21fa2bb9e1500/e2e/b78f0758-4fa2-4641-b5e7-7bb1f6d160c8.md"; F2Url = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ca43f8f47072b8714cbc39ec6f2cd330c0409f4b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b78f0758-4fa2-4641-b5e7-7bb1f6d160c8.cc3665e5c12c629e1fd08d7de85a74ca237eb321.de-de.xlf"; A3Url = "https://github.com/OpenLocalizationTest/oltest/blob/efb688dfdd2528f758c385524c64ff5590bae1c0/e2e/56c59f66-ca12-44db-805b-c48245e4767b.md"; C3Url = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6eb3442b2c767a51f9e30ff60f4edfce947908e9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/56c59f66-ca12-44db-805b-c48245e4767b.fbc4bafe5dd05374473974c13a876717daac1cad.de-de.xlf" }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # push the ".localization-config" row from row 4 down to row 6
    $ws.Cells.Item(6,1).Value = $ws.Cells.Item(4,1).Value()
    $ws.Cells.Item(6,2).Value = $ws.Cells.Item(4,2).Value()
    $ws.Cells.Item(6,4).Value = $ws.Cells.Item(4,4).Value()
    $ws.Cells.Item(6,7).Value = $ws.Cells.Item(4,7).Value()
    $ws.Cells.Item(6,8).Value = $ws.Cells.Item(4,8).Value()

    # new row for file1 (row 4)
    $ws.Cells.Item(4,1).Value = $file1 + ".md"
    $ws.Cells.Item(4,2).Value = $statusReady
    $ws.Cells.Item(4,3).Value = $file1 + "." + $hash1 + "." + $lang.Name + ".xlf"
    $ws.Cells.Item(4,4).Value = $lang.HandoffDT
    $ws.Cells.Item(4,7).Value = "0001-01-01 00:00:00"
    $ws.Cells.Item(4,8).Value = "Include"

    # new row for file2 (row 5)
    $ws.Cells.Item(5,1).Value = $file2 + ".md"
    $ws.Cells.Item(5,2).Value = $statusReady
    $ws.Cells.Item(5,3).Value = $file2 + "." + $hash2 + "." + $lang.Name + ".xlf"
    $ws.Cells.Item(5,4).Value = $lang.HandoffDT
    $ws.Cells.Item(5,7).Value = "0001-01-01 00:00:00"
    $ws.Cells.Item(5,8).Value = "Include"

    # rebuild hyperlinks in final order: A2,C2,E2,F2,A3,C3,A4,C4,A5,C5,A6
    $ws.Range("A1").Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $lang.A2Url, "", "", "b78f0758-4fa2-4641-b5e7-7bb1f6d160c8.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), $lang.C2Url, "", "", "b78f0758-4fa2-4641-b5e7-7bb1f6d160c8.cc3665e5c12c629e1fd08d7de85a74ca237eb321." + $lang.Name + ".xlf") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("E2"), $lang.E2Url, "", "", "b78f0758-4fa2-4641-b5e7-7bb1f6d160c8.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $lang.F2Url, "", "", "b78f0758-4fa2-4641-b5e7-7bb1f6d160c8.cc3665e5c12c629e1fd08d7de85a74ca237eb321." + $lang.Name + ".xlf") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $lang.A3Url, "", "", "56c59f66-ca12-44db-805b-c48245e4767b.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), $lang.C3Url, "", "", "56c59f66-ca12-44db-805b-c48245e4767b.fbc4bafe5dd05374473974c13a876717daac1cad." + $lang.Name + ".xlf") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A4"), $lang.AUrl1, "", "", $file1 + ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C4"), $lang.CUrl1, "", "", $file1 + "." + $hash1 + "." + $lang.Name + ".xlf") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A5"), $lang.AUrl2, "", "", $file2 + ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C5"), $lang.CUrl2, "", "", $file2 + "." + $hash2 + "." + $lang.Name + ".xlf") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/813411781d43403d6b153b57bceef52504c58b28/.localization-config", "", "", $cfgName) | Out-Null
}

Write-Output "Done"
